$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 993.7143
$ws.Range("J29").Value = 1051.5
$ws.Range("L29").Value = 3154.5
$ws.Range("N29").Value = -3716.5

$ws.Range("H46").Value = 1500
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 1600
$ws.Range("K46").Value = 3600
$ws.Range("L46").Value = 4800
$ws.Range("M46").Value = -3481
$ws.Range("N46").Value = -5038

$ws.Range("H60").Value = 1500
$ws.Range("I60").Value = 1200
$ws.Range("J60").Value = 1600
$ws.Range("K60").Value = 3600
$ws.Range("L60").Value = 4800
$ws.Range("M60").Value = -3116
$ws.Range("N60").Value = -5768

$ws.Range("H64").Value = 3552.3076
$ws.Range("J64").Value = 3768.5715
$ws.Range("L64").Value = 3768.5715
$ws.Range("N64").Value = -4264.5715

$ws.Range("H67").Value = 3552.3076
$ws.Range("J67").Value = 3768.5715
$ws.Range("L67").Value = 3768.5715
$ws.Range("N67").Value = -5484.5715

$ws.Range("H93").Value = 22419.387
$ws.Range("J93").Value = 22419.387
$ws.Range("L93").Value = 22419.387
$ws.Range("N93").Value = -27411.387

$ws.Range("H94").Value = 4024.2
$ws.Range("I94").Value = 4024.2
$ws.Range("K94").Value = 4024.2
$ws.Range("M94").Value = -3573.2

$ws.Range("H98").Value = 7977.0435
$ws.Range("J98").Value = 8704.6
$ws.Range("L98").Value = 8704.6
$ws.Range("N98").Value = -11700.6

$ws.Range("H122").Value = 7977.0435
$ws.Range("J122").Value = 8704.6
$ws.Range("L122").Value = 26113.8
$ws.Range("N122").Value = -31013.8

$ws.Range("H131").Value = 3951.1765
$ws.Range("I131").Value = 3327.5
$ws.Range("J131").Value = 4842.143
$ws.Range("K131").Value = 9982.5
$ws.Range("L131").Value = 14526.429
$ws.Range("M131").Value = -4942.5
$ws.Range("N131").Value = -24606.429

$ws.Range("H137").Value = 1589570.9
$ws.Range("I137").Value = 1833043.4
$ws.Range("J137").Value = 7000
$ws.Range("K137").Value = 5499130.199999999
$ws.Range("L137").Value = 21000
$ws.Range("M137").Value = -5496580.199999999
$ws.Range("N137").Value = -26100

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 28181.857
$ws.Range("J121").Value = 28181.857
$ws.Range("L121").Value = 28181.857
$ws.Range("N121").Value = -31675.857

$ws.Range("H137").Value = 45775
$ws.Range("J137").Value = 45775
$ws.Range("L137").Value = 45775
$ws.Range("N137").Value = -55975

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1850.3846
$ws.Range("I107").Value = 1691.7142
$ws.Range("J107").Value = 2035.5
$ws.Range("K107").Value = 1691.7142
$ws.Range("L107").Value = 2035.5
$ws.Range("M107").Value = 228.2858000000001
$ws.Range("N107").Value = -5875.5

$ws.Range("H137").Value = 49351.6
$ws.Range("J137").Value = 49351.6
$ws.Range("L137").Value = 49351.6
$ws.Range("N137").Value = -59551.6

$ws.Range("H138").Value = 40820.383
$ws.Range("J138").Value = 40820.383
$ws.Range("L138").Value = 40820.383
$ws.Range("N138").Value = -51100.383

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2172.25
$ws.Range("I105").Value = 1880
$ws.Range("K105").Value = 1880
$ws.Range("M105").Value = -133

$ws.Range("H107").Value = 1425
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920

$ws.Range("H112").Value = 30357.684
$ws.Range("J112").Value = 30357.684
$ws.Range("L112").Value = 30357.684
$ws.Range("N112").Value = -33311.684

$ws.Range("H123").Value = 41886
$ws.Range("J123").Value = 41886
$ws.Range("L123").Value = 41886
$ws.Range("N123").Value = -51686

$ws.Range("H139").Value = 85011
$ws.Range("J139").Value = 85011
$ws.Range("L139").Value = 85011
$ws.Range("N139").Value = -95291

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2646
$ws.Range("I49").Value = 2525.75
$ws.Range("J49").Value = 2966.6667
$ws.Range("K49").Value = 7577.25
$ws.Range("L49").Value = 8900.000100000001
$ws.Range("M49").Value = -7421.25
$ws.Range("N49").Value = -9212.000100000001

$ws.Range("H58").Value = 3583.3333
$ws.Range("J58").Value = 5000
$ws.Range("L58").Value = 15000
$ws.Range("N58").Value = -15256

$ws.Range("H101").Value = 4000
$ws.Range("J101").Value = 4000
$ws.Range("L101").Value = 12000
$ws.Range("N101").Value = -16868

$ws.Range("H113").Value = 6250597.5
$ws.Range("I113").Value = 641.4
$ws.Range("K113").Value = 1924.2
$ws.Range("M113").Value = 245.8000000000002

$ws.Range("H131").Value = 808.04
$ws.Range("I131").Value = 312.5
$ws.Range("J131").Value = 828.6875
$ws.Range("K131").Value = 937.5
$ws.Range("L131").Value = 2486.0625
$ws.Range("M131").Value = 4102.5
$ws.Range("N131").Value = -12566.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 8334
$ws.Range("I41").Value = 1493.1428
$ws.Range("J41").Value = 20305.5
$ws.Range("K41").Value = 1493.1428
$ws.Range("L41").Value = 20305.5
$ws.Range("M41").Value = -1138.1428
$ws.Range("N41").Value = -21015.5

$ws.Range("H46").Value = 32248.572
$ws.Range("J46").Value = 33415
$ws.Range("L46").Value = 33415
$ws.Range("N46").Value = -33727

$ws.Range("H137").Value = 40186
$ws.Range("J137").Value = 40186
$ws.Range("L137").Value = 40186
$ws.Range("N137").Value = -50386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7541.227
$ws.Range("I132").Value = 3200.6924
$ws.Range("K132").Value = 9602.0772
$ws.Range("M132").Value = -7072.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2000
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 4000
$ws.Range("M81").Value = -2939

$ws.Range("H84").Value = 2000
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 20000
$ws.Range("M84").Value = -14696

$ws.Range("H109").Value = 28377
$ws.Range("J109").Value = 28377
$ws.Range("L109").Value = 28377
$ws.Range("N109").Value = -31151

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 39800
$ws.Range("J112").Value = 39800
$ws.Range("L112").Value = 39800
$ws.Range("N112").Value = -42754

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4159.3076
$ws.Range("I126").Value = 2794.5
$ws.Range("J126").Value = 4765.8887
$ws.Range("K126").Value = 8383.5
$ws.Range("L126").Value = 14297.6661
$ws.Range("M126").Value = -5913.5
$ws.Range("N126").Value = -19237.6661

Write-Output "Applied profit-table updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets."
